# Auto-generated: apply scheduled-runner market-price refresh to the FFXIV leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 509.44446
$ws.Range("J4").Value = 798.5
$ws.Range("L4").Value = 798.5
$ws.Range("N4").Value = -1026.5

$ws.Range("H80").Value = 548.62964
$ws.Range("I80").Value = 187.76471
$ws.Range("J80").Value = 1162.1
$ws.Range("K80").Value = 563.29413
$ws.Range("L80").Value = 3486.3
$ws.Range("M80").Value = 434.70587
$ws.Range("N80").Value = -5482.299999999999

$ws.Range("H83").Value = 548.62964
$ws.Range("I83").Value = 187.76471
$ws.Range("J83").Value = 1162.1
$ws.Range("K83").Value = 1689.88239
$ws.Range("L83").Value = 10458.9
$ws.Range("M83").Value = 3302.11761
$ws.Range("N83").Value = -20442.9

$ws.Range("H100").Value = 2046.9584
$ws.Range("I100").Value = 1106.909
$ws.Range("K100").Value = 1106.909
$ws.Range("M100").Value = -565.9090000000001

$ws.Range("H115").Value = 6738447
$ws.Range("I115").Value = 7159506
$ws.Range("K115").Value = 21478518
$ws.Range("M115").Value = -21476951

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 188.88889
$ws.Range("I10").Value = 188.88889
$ws.Range("K10").Value = 188.88889
$ws.Range("M10").Value = -18.88889

$ws.Range("H32").Value = 33082.74
$ws.Range("I32").Value = 37660.91
$ws.Range("K32").Value = 37660.91
$ws.Range("M32").Value = -37373.91

$ws.Range("H61").Value = 8133837.5
$ws.Range("I61").Value = 9526430
$ws.Range("K61").Value = 9526430
$ws.Range("M61").Value = -9526218

$ws.Range("H74").Value = 253050.1
$ws.Range("I74").Value = 304479.97
$ws.Range("K74").Value = 304479.97
$ws.Range("M74").Value = -303605.97

$ws.Range("H77").Value = 253050.1
$ws.Range("I77").Value = 304479.97
$ws.Range("K77").Value = 1522399.85
$ws.Range("M77").Value = -1518031.85

$ws.Range("H136").Value = 8133837.5
$ws.Range("I136").Value = 9526430
$ws.Range("K136").Value = 28579290
$ws.Range("M136").Value = -28576740

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3121.625
$ws.Range("I22").Value = 1828.8334
$ws.Range("K22").Value = 1828.8334
$ws.Range("M22").Value = -1655.8334

$ws.Range("H105").Value = 29420612
$ws.Range("I105").Value = 37047470
$ws.Range("J105").Value = 2723.8572
$ws.Range("K105").Value = 37047470
$ws.Range("L105").Value = 2723.8572
$ws.Range("M105").Value = -37045723
$ws.Range("N105").Value = -6217.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 23333
$ws.Range("J3").Value = 22500
$ws.Range("L3").Value = 22500
$ws.Range("N3").Value = -22726

$ws.Range("H62").Value = 11437.692
$ws.Range("I62").Value = 8278.666999999999
$ws.Range("J62").Value = 14145.429
$ws.Range("K62").Value = 8278.666999999999
$ws.Range("L62").Value = 14145.429
$ws.Range("M62").Value = -7654.666999999999
$ws.Range("N62").Value = -15393.429

$ws.Range("H65").Value = 11437.692
$ws.Range("I65").Value = 8278.666999999999
$ws.Range("J65").Value = 14145.429
$ws.Range("K65").Value = 41393.335
$ws.Range("L65").Value = 70727.145
$ws.Range("M65").Value = -38273.335
$ws.Range("N65").Value = -76967.145

$ws.Range("H125").Value = 59499.5
$ws.Range("J125").Value = 59499.5
$ws.Range("L125").Value = 59499.5
$ws.Range("N125").Value = -64419.5

$ws.Range("H132").Value = 19194.352
$ws.Range("I132").Value = 3929.1538
$ws.Range("K132").Value = 11787.4614
$ws.Range("M132").Value = -9257.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3482.25
$ws.Range("I2").Value = 6.1
$ws.Range("J2").Value = 6958.4
$ws.Range("K2").Value = 36.59999999999999
$ws.Range("L2").Value = 41750.39999999999
$ws.Range("M2").Value = 76.40000000000001
$ws.Range("N2").Value = -41976.39999999999

$ws.Range("H7").Value = 339.2143
$ws.Range("I7").Value = 299.8
$ws.Range("J7").Value = 361.1111
$ws.Range("K7").Value = 899.4000000000001
$ws.Range("L7").Value = 1083.3333
$ws.Range("M7").Value = -787.4000000000001
$ws.Range("N7").Value = -1307.3333

$ws.Range("H34").Value = 1204.5834
$ws.Range("I34").Value = 116.416664
$ws.Range("J34").Value = 2292.75
$ws.Range("K34").Value = 349.249992
$ws.Range("L34").Value = 6878.25
$ws.Range("M34").Value = -265.249992
$ws.Range("N34").Value = -7046.25

$ws.Range("H38").Value = 26.2
$ws.Range("I38").Value = 36
$ws.Range("J38").Value = 18.181818
$ws.Range("K38").Value = 108
$ws.Range("L38").Value = 54.545454
$ws.Range("M38").Value = 239
$ws.Range("N38").Value = -748.5454539999999

$ws.Range("H39").Value = 2887.9333
$ws.Range("I39").Value = 478.57144
$ws.Range("J39").Value = 4996.125
$ws.Range("K39").Value = 1435.71432
$ws.Range("L39").Value = 14988.375
$ws.Range("M39").Value = -1141.71432
$ws.Range("N39").Value = -15576.375

$ws.Range("H55").Value = 1409
$ws.Range("I55").Value = 1456.3334
$ws.Range("J55").Value = 1380.6
$ws.Range("K55").Value = 4369.0002
$ws.Range("L55").Value = 4141.799999999999
$ws.Range("M55").Value = -4192.0002
$ws.Range("N55").Value = -4495.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 28814.572
$ws.Range("I2").Value = 307.88
$ws.Range("K2").Value = 307.88
$ws.Range("M2").Value = -194.88

$ws.Range("H9").Value = 4399.3335
$ws.Range("I9").Value = 199
$ws.Range("J9").Value = 6499.5
$ws.Range("K9").Value = 199
$ws.Range("L9").Value = 6499.5
$ws.Range("M9").Value = -29
$ws.Range("N9").Value = -6839.5

$ws.Range("H43").Value = 3999.8333
$ws.Range("I43").Value = 800
$ws.Range("J43").Value = 19999
$ws.Range("K43").Value = 800
$ws.Range("L43").Value = 19999
$ws.Range("M43").Value = -649
$ws.Range("N43").Value = -20301

$ws.Range("H92").Value = 15083.667
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 15083.667
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 15083.667
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -18827.667

$ws.Range("H123").Value = 42479.125
$ws.Range("J123").Value = 59833
$ws.Range("L123").Value = 59833
$ws.Range("N123").Value = -64733

$ws.Range("H128").Value = 48333
$ws.Range("J128").Value = 48333
$ws.Range("L128").Value = 48333
$ws.Range("N128").Value = -58293

$ws.Range("H134").Value = 101219.8
$ws.Range("J134").Value = 101219.8
$ws.Range("L134").Value = 303659.4
$ws.Range("N134").Value = -308729.4

$ws.Range("H136").Value = 43055.938
$ws.Range("J136").Value = 43055.938
$ws.Range("L136").Value = 129167.814
$ws.Range("N136").Value = -134267.814

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 5298.3335
$ws.Range("I10").Value = 997.5
$ws.Range("J10").Value = 13900
$ws.Range("K10").Value = 997.5
$ws.Range("L10").Value = 13900
$ws.Range("M10").Value = -857.5
$ws.Range("N10").Value = -14180

$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H20").Value = 7501000
$ws.Range("J20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("N20").Value = -2452

$ws.Range("H22").Value = 3172.9
$ws.Range("I22").Value = 1862.4584
$ws.Range("K22").Value = 1862.4584
$ws.Range("M22").Value = -1567.4584

$ws.Range("H27").Value = 3172.9
$ws.Range("I27").Value = 1862.4584
$ws.Range("K27").Value = 1862.4584
$ws.Range("M27").Value = -1755.4584

$ws.Range("H46").Value = 8054.5
$ws.Range("I46").Value = 1300
$ws.Range("K46").Value = 1300
$ws.Range("M46").Value = -1112

$ws.Range("H68").Value = 4399.1
$ws.Range("I68").Value = 3081.8333
$ws.Range("K68").Value = 3081.8333
$ws.Range("M68").Value = -2332.8333

$ws.Range("H71").Value = 4399.1
$ws.Range("I71").Value = 3081.8333
$ws.Range("K71").Value = 15409.1665
$ws.Range("M71").Value = -11665.1665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10690.525
$ws.Range("I81").Value = 4682.5625
$ws.Range("J81").Value = 14695.833
$ws.Range("K81").Value = 9365.125
$ws.Range("L81").Value = 29391.666
$ws.Range("M81").Value = -8304.125
$ws.Range("N81").Value = -31513.666

$ws.Range("H84").Value = 10690.525
$ws.Range("I84").Value = 4682.5625
$ws.Range("J84").Value = 14695.833
$ws.Range("K84").Value = 46825.625
$ws.Range("L84").Value = 146958.33
$ws.Range("M84").Value = -41521.625
$ws.Range("N84").Value = -157566.33

$ws.Range("H100").Value = 2530.6155
$ws.Range("I100").Value = 2199.3333
$ws.Range("K100").Value = 4398.6666
$ws.Range("M100").Value = -3857.6666

$ws.Range("H136").Value = 3486235.8
$ws.Range("I136").Value = 3862318
$ws.Range("K136").Value = 11586954
$ws.Range("M136").Value = -11584404
